$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2025-03-09 06:05:12"
$ws.Range("F2").Value = "http://49.234.6.241:5230/api/v1/resource/16"
$ws.Range("M2").Value = 0.004
$ws.Range("N2").Value = 1
$ws.Range("Q2").Value = $false

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("D3").Value = "2025-03-09 06:05:12"
$ws.Range("F3").Value = "http://49.234.6.241:5230/api/v1/resource/16"
$ws.Range("G3").Value = "/api/v1/resource/16"
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = $false

# Row 4
$ws.Range("D4").Value = "2025-03-09 06:05:12"
$ws.Range("M4").Value = 0.003
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $true
$ws.Range("Q4").Value = $true

# Row 5
$ws.Range("D5").Value = "2025-03-09 06:05:12"
$ws.Range("F5").Value = "http://47.97.114.24:5230/api/v1/resource/16"
$ws.Range("G5").Value = "/api/v1/resource/16"
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $true
$ws.Range("Q5").Value = $true

# Row 6
$ws.Range("D6").Value = "2025-03-09 06:05:12"
$ws.Range("F6").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G6").Value = "/api/v1/memo/21"
$ws.Range("M6").Value = 0.003
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = $false
